$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.841.63'
$ws.Range('E2').Value = '  +2.10%  '

$ws.Range('D3').Value = '3.473.38'
$ws.Range('E3').Value = '  +2.48%  '

$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = '581.25'
$ws.Range('E5').Value = '  +0.97%  '

$ws.Range('D6').Value = '147.19'
$ws.Range('E6').Value = '  +4.48%  '

$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('E8').Value = '  +1.35%  '

$ws.Range('E9').Value = '  -1.14%  '

$ws.Range('D10').Value = '0.125'
$ws.Range('E10').Value = '  +2.12%  '

$ws.Range('E11').Value = '  +3.73%  '

$ws.Range('D12').Value = '4.069.06'
$ws.Range('E12').Value = '  +2.47%  '

$ws.Range('D13').Value = '29.92'
$ws.Range('E13').Value = '  +4.98%  '

$ws.Range('D14').Value = '0.127'
$ws.Range('E14').Value = '  +1.30%  '

$ws.Range('D15').Value = '3.473.02'
$ws.Range('E15').Value = '  +2.15%  '

$ws.Range('E16').Value = '  +0.96%  '

$ws.Range('D17').Value = '62.892.18'
$ws.Range('E17').Value = '  +2.17%  '

$ws.Range('D18').Value = '6.35'
$ws.Range('E18').Value = '  +3.46%  '

$ws.Range('D19').Value = '14.46'
$ws.Range('E19').Value = '  +5.94%  '

$ws.Range('D20').Value = '9.29'
$ws.Range('E20').Value = '  +3.20%  '

$ws.Range('D21').Value = '390.07'
$ws.Range('E21').Value = '  -0.26%  '

$ws.Range('D22').Value = '0.566'
$ws.Range('E22').Value = '  +2.27%  '

$ws.Range('D23').Value = '74.94'
$ws.Range('E23').Value = '  -0.44%  '

$ws.Range('E24').Value = '  -0.08%  '

$ws.Range('D25').Value = '3.619.54'
$ws.Range('E25').Value = '  +2.54%  '

$ws.Range('D26').Value = '0.0000115'
$ws.Range('E26').Value = '  +2.50%  '

$ws.Range('D27').Value = '0.181'
$ws.Range('E27').Value = '  -6.95%  '

$ws.Range('D28').Value = '7.67'
$ws.Range('E28').Value = '  +5.73%  '

$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.73%  '

$ws.Range('D30').Value = '8.22'
$ws.Range('E30').Value = '  +2.07%  '

$ws.Range('E31').Value = '  +0.55%  '

$ws.Range('D32').Value = '1.40'
$ws.Range('E32').Value = '  -0.33%  '

$ws.Range('E33').Value = '  -0.01%  '

$ws.Range('D34').Value = '23.80'
$ws.Range('E34').Value = '  +2.24%  '

$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').Value = '7.11'
$ws.Range('E35').Value = '  +2.69%  '

$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').Value = '5.29'
$ws.Range('E36').Value = '  +4.69%  '

$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '171.20'
$ws.Range('E37').Value = '  +1.93%  '

$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D38').Value = '31.39'
$ws.Range('E38').Value = '  +19.72%  '

$ws.Range('D39').Value = '1.56'
$ws.Range('E39').Value = '  +6.38%  '

$ws.Range('D40').Value = '3.512.05'
$ws.Range('E40').Value = '  +2.63%  '

$ws.Range('D41').Value = '0.0770'
$ws.Range('E41').Value = '  +0.20%  '

$ws.Range('D42').Value = '0.801'
$ws.Range('E42').Value = '  +2.84%  '

$ws.Range('D43').Value = '4.51'
$ws.Range('E43').Value = '  +1.95%  '

$ws.Range('D44').Value = '42.20'
$ws.Range('E44').Value = '  -0.77%  '

$ws.Range('D45').Value = '1.71'
$ws.Range('E45').Value = '  +3.49%  '

$ws.Range('E46').Value = '  +3.09%  '

$ws.Range('D47').Value = '2.602.00'
$ws.Range('E47').Value = '  +6.05%  '

$ws.Range('D48').Value = '23.63'
$ws.Range('E48').Value = '  +2.59%  '

$ws.Range('D49').Value = '2.27'
$ws.Range('E49').Value = '  +11.88%  '

$ws.Range('D50').Value = '6.79'
$ws.Range('E50').Value = '  +1.80%  '

$ws.Range('B51').Value = 'FirstDigitalUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  +0.10%  '
